# Remove the decorative horizontal-rule paragraphs (VML <v:rect ... o:hr="t"/>
# inside a <w:pict>) that separate each resume section. These paragraphs carry
# no text of their own -- they only contain the picture run -- so we find them
# by inspecting each paragraph's underlying OOXML for the horizontal-rule
# marker and delete the whole paragraph (including its paragraph mark).

$d = $word.ActiveDocument

$targets = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    $xml = $p.Range.WordOpenXML
    if ($xml -like "*o:hr=*") {
        [void]$targets.Add($p)
    }
}

for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $targets[$i].Range.Delete()
}
